$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header date (D1) merged like B1:C1, value 41728 (2014-03-30) ---
$ws.Range("D1:E1").Merge()
$ws.Range("D1").Value = 41728
$ws.Range("D1:E1").HorizontalAlignment = -4108
$ws.Range("D1").NumberFormat = "d-mmm"

# --- Correct B3 (Control plant leaf count) from 19 to 18 ---
$ws.Range("B3").Value = 18

# --- New experimental-plant leaf counts in column D ---
$ws.Range("D3").Value = 18
$ws.Range("D4").Value = 12
$ws.Range("D5").Value = 10
$ws.Range("D6").Value = 20

# --- New notes in column E ---
$ws.Range("E4").Value = "Leaves on tip little wrinkly"
$ws.Range("E5").Value = "Leaves on tip little wrinkly"
$ws.Range("E6").Value = "Leaves on tip more wrinkly"

# --- Column E width to match new Notes column (~24.14 chars; engine quantizes
#     ColumnWidth to 1/6-character pixel steps, so 23.25 is the closest input
#     that lands on the nearest achievable stored width) ---
$ws.Range("E1").ColumnWidth = 23.25

# --- Selection left where the edit finished ---
$ws.Range("E2").Select()
